$wb = $excel.ActiveWorkbook

# --- Sheet "Sampling Events": fix typo in eventID (B2) ---
$wsEvents = $wb.Worksheets.Item("Sampling Events")
$wsEvents.Range("B2").Value = "UNIPA-2002AS-HS007-NA001"

# --- Sheet "Occurrences": replace the single aggregated record with
#     8 individual occurrence records (minus 3 transcriptions that used
#     to be duplicated under the old KM001 event) ---
$wsOcc = $wb.Worksheets.Item("Occurrences")

$eventID = "UNIPA-2002AS-HS007-NA001"
$basisOfRecord = "Human Observation"
$eventDate = "2001-12-27/2002-02-25"
$kingdom = "Animalia"
$scientificName = "Spilocuscus maculatus"
$taxonRank = "Spesies"
$vernacularName = "Rambab magnan/Rambab mangawak"
$decimalLatitude = "'-1.076681"
$decimalLongitude = "'134.834596"
$geodeticDatum = "WGS84"
$countryCode = "ID"
$individualCount = 1

for ($i = 1; $i -le 8; $i++) {
    $row = $i + 1
    $occurrenceID = "UNIPA-2002AS-HS007-NA001-VE00$i"

    $wsOcc.Range("A$row").Value = $eventID
    $wsOcc.Range("B$row").Value = $occurrenceID
    $wsOcc.Range("C$row").Value = $basisOfRecord
    $wsOcc.Range("D$row").Value = $eventDate
    $wsOcc.Range("E$row").Value = $kingdom
    $wsOcc.Range("F$row").Value = $scientificName
    $wsOcc.Range("G$row").Value = $taxonRank
    $wsOcc.Range("H$row").Value = $vernacularName
    $wsOcc.Range("I$row").Value = $decimalLatitude
    $wsOcc.Range("J$row").Value = $decimalLongitude
    $wsOcc.Range("K$row").Value = $geodeticDatum
    $wsOcc.Range("L$row").Value = $countryCode
    $wsOcc.Range("M$row").Value = $individualCount
}

# Row 2 additionally kept the occurrenceStatus / remarks notes that were
# originally attached to the single aggregated record.
$wsOcc.Range("P2").Value = "2001-10/2001-11"
$wsOcc.Range("Q2").Value = "Didalam Skripsi tidak disebutkan tanggal penelitian"

# Clear the remark/status values that used to sit on the old single row
# for rows 3-9 (they never had them before either, so nothing to do).
